$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "department" column (C) value for the course row was simplified
# from the old long-form label to the shorter value used going forward.
$ws.Range("C2").Value = "English"
